$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) values per row to reflect the refreshed
# cryptos snapshot. Numeric-looking Price strings are forced back to text
# (NumberFormat "@") so Excel does not auto-convert them to real numbers,
# then the style is reset to "Normal" so no stray cell formatting is left behind.

$ws.Range("D2").Value = '58.948.55'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").Value = '2.502.71'
$ws.Range("E3").Value = '  -0.72%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("E8").Value = '  +1.06%  '

$ws.Range("D9").Value = '2.528.10'
$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("E10").Value = '  +1.19%  '

$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.37%  '

$ws.Range("E13").Value = '  -2.08%  '

$ws.Range("D14").Value = '2.961.94'
$ws.Range("E14").Value = '  -0.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.89%  '

$ws.Range("D16").Value = '58.929.36'
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").Value = '2.515.44'
$ws.Range("E18").Value = '  -1.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.92'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.49%  '

$ws.Range("E25").Value = '  -0.68%  '

$ws.Range("E26").Value = '  +0.62%  '

$ws.Range("E27").Value = '  +1.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.41%  '

$ws.Range("D29").Value = '0.0₃0778'
$ws.Range("E29").Value = '  +1.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.26%  '

$ws.Range("E31").Value = '  -1.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '168.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.57%  '

$ws.Range("E33").Value = '  +6.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.66%  '

$ws.Range("E37").Value = '  -1.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.835'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.54%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '283.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("E44").Value = '  -0.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("E49").Value = '  +0.39%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0224'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.56%  '

$ws.Range("E51").Value = '  +0.17%  '

# Rows 45 and 46 swap places in ranking: Mantle <-> Aave
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.00%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.607"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.79%  "
